$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ShipmentTrackNum / PackageTrackNum values (text, not numbers) that
# replace the stale ones in column C (rows 2-22), mirrored in column D for
# the rows where PackageTrackNum == ShipmentTrackNum. Ordered top to bottom
# to match the row layout.
$rows      = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)
$trackNums = @(
    "320018543178",
    "320018543226",
    "320018543259",
    "320018543270",
    "320018543318",
    "320018543330",
    "320018543362",
    "320018543384",
    "320018543410",
    "320018543432",
    "320018543476",
    "320018543498",
    "320018543524",
    "320018543546",
    "320018543579",
    "320018548011",
    "320018548055",
    "320018548077",
    "320018548103",
    "320018548125",
    "320018548158"
)

# Rows where the PackageTrackNum (column D) mirrors the ShipmentTrackNum
# (column C) value.
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row   = $rows[$i]
    $value = $trackNums[$i]

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $value
    $cCell.Style = "Normal"

    if ($mirrorRows -contains $row) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $value
        $dCell.Style = "Normal"
    }
}

# Q3 flips from FAIL to PASS.
$ws.Range("Q3").Value = "PASS"
